$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Apply formatting first by copying from existing same-column cells so the
#     same cellXf / font indices are reused instead of new ones being minted ---
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

$ws.Range("I16").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$ws.Range("I4").PasteSpecial(-4122)
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I14").PasteSpecial(-4122)

# --- New header cell I1 ---
$ws.Range("I1").Value = "new Response"

# --- New response JSON cells in column I ---
$ws.Range("I3").Value = "[`n    {`n        ""city"": ""San Jose"",`n        ""code"": ""SJC"",`n        ""id"": ""61849d3f4367d925b16ff24b"",`n        ""name"": ""San Jose International Airport""`n    },`n    {`n        ""city"": ""San Francisco"",`n        ""code"": ""SFO"",`n        ""id"": ""61849d5f4367d925b16ff24c"",`n        ""name"": ""San Francisco International Airport""`n    }]"

$ws.Range("I4").Value = "{`n    ""city"": ""San Jose"",`n    ""code"": ""SJC"",`n    ""id"": ""61849d3f4367d925b16ff24b"",`n    ""name"": ""San Jose International Airport""`n}"

$ws.Range("I13").Value = "[`n    {`n        ""id"": ""61a5aac8b657dbd44f933bb9"",`n        ""name"": ""Airbus A320"",`n        ""seat_chart"": {`n            ""aisle"": [`n                ""1A""`n            ],`n            ""middle"": [`n                ""1B""`n            ],`n            ""window"": [`n                ""1C""`n            ]`n        },`n        ""total_seats"": 90`n    },`n    {`n        ""id"": ""61a5ab60b657dbd44f933bba"",`n        ""name"": ""Boeing 747"",`n        ""seat_chart"": {`n            ""aisle"": [`n                ""1A""`n            ],`n            ""middle"": [`n                ""1B""`n            ],`n            ""window"": [`n                ""1C""`n            ]`n        },`n        ""total_seats"": 90`n    },`n    {`n        ""id"": ""61a5abf6b657dbd44f933bbb"",`n        ""name"": ""Airbus A520"",`n        ""seat_chart"": {`n            ""aisle"": [`n                ""1A""`n            ],`n            ""middle"": [`n                ""1B""`n            ],`n            ""window"": [`n                ""1C""`n            ]`n        },`n        ""total_seats"": 90`n    }`n]"

$ws.Range("I14").Value = "{`n    ""id"": ""61a5aac8b657dbd44f933bb9"",`n    ""name"": ""Airbus A320"",`n    ""seat_chart"": {`n        ""aisle"": [`n            ""1A""`n        ],`n        ""middle"": [`n            ""1B""`n        ],`n        ""window"": [`n            ""1C""`n        ]`n    },`n    ""total_seats"": 90`n}"

# --- Row heights grow to fit the new wrapped JSON text in column I ---
$ws.Rows.Item(3).RowHeight = 295
$ws.Rows.Item(13).RowHeight = 409.6
$ws.Rows.Item(14).RowHeight = 192

# --- Column widths: H shrinks, I gets an explicit width ---
# (target stored widths are 11.1640625 / 24.33203125; the engine quantises
#  ColumnWidth to 1/6-character steps on write, so these inputs land on the
#  closest reachable stored width: 11.1666.. / 24.3333..)
$ws.Columns.Item(8).ColumnWidth = 10.333333333333334
$ws.Columns.Item(9).ColumnWidth = 23.5

# --- View: scroll position + active selection ---
$ws.Range("A13").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I14").Select()
